$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 17.470401
$ws.Range("H2").Value = 52.411203
$ws.Range("I2").Value = 0.8600988665959021
$ws.Range("J2").Value = 0.8884442399952684
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 44.13164066666667
$ws.Range("N2").Value = 132.394922
$ws.Range("O2").Value = 0.4415399811720331
$ws.Range("P2").Value = 0.4562856844211927
$ws.Range("Q2").Value = 770.997459234574
$ws.Range("R2").Value = 6938.977133111166
$ws.Range("S2").Value = 0.3797680373628416
$ws.Range("T2").Value = 0.4053843881163074
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 17.470401
$ws.Range("H3").Value = 52.411203
$ws.Range("I3").Value = 0.8600988665959021
$ws.Range("J3").Value = 0.8884442399952684
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 14.93259333333333
$ws.Range("N3").Value = 44.79778
$ws.Range("O3").Value = 0.1494015830739255
$ws.Range("P3").Value = 0.1543910098595022
$ws.Range("Q3").Value = 260.87839350326
$ws.Range("R3").Value = 2347.90554152934
$ws.Range("S3").Value = 0.1285001322695168
$ws.Range("T3").Value = 0.1371678034167274
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 17.470401
$ws.Range("H4").Value = 52.411203
$ws.Range("I4").Value = 0.8600988665959021
$ws.Range("J4").Value = 0.8884442399952684
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.61024133333333
$ws.Range("N4").Value = 40.830724
$ws.Range("O4").Value = 0.1361713639304118
$ws.Range("P4").Value = 0.1407189532975654
$ws.Range("Q4").Value = 237.776373800108
$ws.Range("R4").Value = 2139.987364200972
$ws.Range("S4").Value = 0.1171208357793653
$ws.Range("T4").Value = 0.1250209435153852
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 17.470401
$ws.Range("H5").Value = 52.411203
$ws.Range("I5").Value = 0.8600988665959021
$ws.Range("J5").Value = 0.8884442399952684
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.584752
$ws.Range("N5").Value = 52.754256
$ws.Range("O5").Value = 0.1759366057935712
$ws.Range("P5").Value = 0.1818121982434553
$ws.Range("Q5").Value = 307.212668925552
$ws.Range("R5").Value = 2764.914020329968
$ws.Range("S5").Value = 0.1513228752357806
$ws.Range("T5").Value = 0.1615300002902758
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.470401
$ws.Range("H6").Value = 52.411203
$ws.Range("I6").Value = 0.8600988665959021
$ws.Range("J6").Value = 0.8884442399952684
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 9.690137500000001
$ws.Range("N6").Value = 19.380275
$ws.Range("O6").Value = 0.09695046603005844
$ws.Range("P6").Value = 0.06679215417828435
$ws.Range("Q6").Value = 169.2905878701375
$ws.Range("R6").Value = 1015.743527220825
$ws.Range("S6").Value = 0.08338698594839777
$ws.Range("T6").Value = 0.05934110465657264
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.8975426666666667
$ws.Range("H7").Value = 2.692628
$ws.Range("I7").Value = 0.04418761940962108
$ws.Range("J7").Value = 0.04564386429080782
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 44.13164066666667
$ws.Range("N7").Value = 132.394922
$ws.Range("O7").Value = 0.4415399811720331
$ws.Range("P7").Value = 0.4562856844211927
$ws.Range("Q7").Value = 39.61003044833512
$ws.Range("R7").Value = 356.490274035016
$ws.Range("S7").Value = 0.01951060064216106
$ws.Range("T7").Value = 0.02082664185755928
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.8975426666666667
$ws.Range("H8").Value = 2.692628
$ws.Range("I8").Value = 0.04418761940962108
$ws.Range("J8").Value = 0.04564386429080782
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.93259333333333
$ws.Range("N8").Value = 44.79778
$ws.Range("O8").Value = 0.1494015830739255
$ws.Range("P8").Value = 0.1543910098595022
$ws.Range("Q8").Value = 13.40263964064889
$ws.Range("R8").Value = 120.62375676584
$ws.Range("S8").Value = 0.006601700292065507
$ws.Range("T8").Value = 0.007047002301747889
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.8975426666666667
$ws.Range("H9").Value = 2.692628
$ws.Range("I9").Value = 0.04418761940962108
$ws.Range("J9").Value = 0.04564386429080782
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.61024133333333
$ws.Range("N9").Value = 40.830724
$ws.Range("O9").Value = 0.1361713639304118
$ws.Range("P9").Value = 0.1407189532975654
$ws.Range("Q9").Value = 12.21577230029689
$ws.Range("R9").Value = 109.941950702672
$ws.Range("S9").Value = 0.006017088403846039
$ws.Range("T9").Value = 0.0064229568074586
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8975426666666667
$ws.Range("H10").Value = 2.692628
$ws.Range("I10").Value = 0.04418761940962108
$ws.Range("J10").Value = 0.04564386429080782
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 17.584752
$ws.Range("N10").Value = 52.754256
$ws.Range("O10").Value = 0.1759366057935712
$ws.Range("P10").Value = 0.1818121982434553
$ws.Range("Q10").Value = 15.783065202752
$ws.Range("R10").Value = 142.047586824768
$ws.Range("S10").Value = 0.007774219777026862
$ws.Range("T10").Value = 0.008298611303037724
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.8975426666666667
$ws.Range("H11").Value = 2.692628
$ws.Range("I11").Value = 0.04418761940962108
$ws.Range("J11").Value = 0.04564386429080782
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 9.690137500000001
$ws.Range("N11").Value = 19.380275
$ws.Range("O11").Value = 0.09695046603005844
$ws.Range("P11").Value = 0.06679215417828435
$ws.Range("Q11").Value = 8.697311852116668
$ws.Range("R11").Value = 52.1838711127
$ws.Range("S11").Value = 0.00428401029452162
$ws.Range("T11").Value = 0.003048652021004324
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 1.944141
$ws.Range("H12").Value = 3.888282
$ws.Range("I12").Value = 0.09571351399447693
$ws.Range("J12").Value = 0.06591189571392365
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 44.13164066666667
$ws.Range("N12").Value = 132.394922
$ws.Range("O12").Value = 0.4415399811720331
$ws.Range("P12").Value = 0.4562856844211927
$ws.Range("Q12").Value = 85.79813201733401
$ws.Range("R12").Value = 514.7887921040041
$ws.Range("S12").Value = 0.04226134316703047
$ws.Range("T12").Value = 0.03007465444732593
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 1.944141
$ws.Range("H13").Value = 3.888282
$ws.Range("I13").Value = 0.09571351399447693
$ws.Range("J13").Value = 0.06591189571392365
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 14.93259333333333
$ws.Range("N13").Value = 44.79778
$ws.Range("O13").Value = 0.1494015830739255
$ws.Range("P13").Value = 0.1543910098595022
$ws.Range("Q13").Value = 29.03106693566
$ws.Range("R13").Value = 174.18640161396
$ws.Range("S13").Value = 0.01429975051234317
$ws.Range("T13").Value = 0.01017620414102686
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 1.944141
$ws.Range("H14").Value = 3.888282
$ws.Range("I14").Value = 0.09571351399447693
$ws.Range("J14").Value = 0.06591189571392365
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 13.61024133333333
$ws.Range("N14").Value = 40.830724
$ws.Range("O14").Value = 0.1361713639304118
$ws.Range("P14").Value = 0.1407189532975654
$ws.Range("Q14").Value = 26.460228196028
$ws.Range("R14").Value = 158.761369176168
$ws.Range("S14").Value = 0.01303343974720048
$ws.Range("T14").Value = 0.009275052974721625
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 1.944141
$ws.Range("H15").Value = 3.888282
$ws.Range("I15").Value = 0.09571351399447693
$ws.Range("J15").Value = 0.06591189571392365
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 17.584752
$ws.Range("N15").Value = 52.754256
$ws.Range("O15").Value = 0.1759366057935712
$ws.Range("P15").Value = 0.1818121982434553
$ws.Range("Q15").Value = 34.187237338032
$ws.Range("R15").Value = 205.123424028192
$ws.Range("S15").Value = 0.01683951078076375
$ws.Range("T15").Value = 0.01198358665014184
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 1.944141
$ws.Range("H16").Value = 3.888282
$ws.Range("I16").Value = 0.09571351399447693
$ws.Range("J16").Value = 0.06591189571392365
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 9.690137500000001
$ws.Range("N16").Value = 19.380275
$ws.Range("O16").Value = 0.09695046603005844
$ws.Range("P16").Value = 0.06679215417828435
$ws.Range("Q16").Value = 18.8389936093875
$ws.Range("R16").Value = 75.35597443755
$ws.Range("S16").Value = 0.009279469787139059
$ws.Range("T16").Value = 0.008298611303037724
